$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day-4 block (rows 14-17), mirroring the existing Day-1/2/3 blocks ---
# Row 14: Day label (top-left of the merged A14:A17) + Shubhendu
$ws.Range("A14").Value = "Day-4"
$ws.Range("B14").Value = "Shubhendu"
$ws.Range("C14").Value = "Rest API"
$ws.Range("D14").Value = "In progress"

# Row 15: Praneetha
$ws.Range("B15").Value = "Praneetha"
$ws.Range("C15").Value = "Rest API"
$ws.Range("D15").Value = "In progress"

# Row 16: Shravya
$ws.Range("B16").Value = "Shravya"
$ws.Range("C16").Value = "Rest API"
$ws.Range("D16").Value = "In progress"

# Row 17: Shafeeq
$ws.Range("B17").Value = "Shafeeq"
$ws.Range("C17").Value = "Rest API"
$ws.Range("D17").Value = "In progress"

# --- Match the B-column style used by the other day blocks (B6, B10, ...) ---
$ws.Range("B14").Style = $ws.Range("B10").Style
$ws.Range("B15").Style = $ws.Range("B11").Style
$ws.Range("B16").Style = $ws.Range("B12").Style
$ws.Range("B17").Style = $ws.Range("B13").Style

# --- New (mostly empty) column E picks up a default width definition ---
$ws.Columns.Item(5).ColumnWidth = 8.333333333333334

# --- Selection moves on to the next empty row, as it would after data entry ---
$ws.Range("D18").Select()
